$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "...si se trata de un botón el nombre..." -> "...un label el nombre..."
#    and move the "_GoBack" bookmark so it sits right after the new
#    "label" word (it currently lives alone in the trailing empty
#    paragraph).
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("botón", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# $r now covers the found word ("botón"); replace its text with "label".
$r.Text = "label"

# Toggling Bold briefly forces the run containing "label" to become its
# own run (distinct from the surrounding, identically-formatted text)
# so the bookmark we are about to insert lands between two run
# boundaries instead of in the middle of a merged run.
$r.Bold = 1
$bmRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
$r.Bold = 0

# ------------------------------------------------------------------
# 2) "botonSalir" -> "labelNombreJugador" (keeps bold formatting and
#    the surrounding spell-check markers intact).
# ------------------------------------------------------------------
$d.Content.Find.Execute("botonSalir", $true, $false, $false, $false, $false, $true, 1, $false, "labelNombreJugador", 2) | Out-Null
